# Apply the "2022 -> 2023" country/pays codes refresh:
#  1. Rename the worksheet (and, automatically, the _FilterDatabase
#     defined name that references it) from the 2022 to the 2023 label.
#  2. Merge the separate English "Turkey" / French "Turquie" shared
#     strings into the single localized name "Türkiye" used for both
#     the English and French columns on the TUR row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename sheet (defined names referencing the sheet update automatically).
$ws.Name = "AfDD_DDAf_2023_CntryPaysCodes"

# 2) Find the TUR row (column A = "TUR") and update the English/French
#    country-name columns (B/C) to the unified "Türkiye" name.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$turRow = 0
for ($r = 1; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 1).Value() -eq "TUR") {
        $turRow = $r
        break
    }
}

if ($turRow -gt 0) {
    $ws.Cells.Item($turRow, 2).Value = "Türkiye"
    $ws.Cells.Item($turRow, 3).Value = "Türkiye"
}
